# "Fixed number of slides"
#
# The title slide (slide 1) and the closing "Thank you" slide (slide 8)
# each had an extra slide-number placeholder that shouldn't have been
# there (those two slides are not meant to show a page number). Turn the
# per-slide slide-number display off for just those two slides, which
# removes the stray placeholder shape from each.
#
# Also update the outline color of the two rounded-rectangle callouts on
# slide 4 from the theme's dk2 scheme color to explicit RGB colors.

$p = $ppt.ActivePresentation

function Set-LineColorRGB($shape, [string]$hex) {
    $r = [Convert]::ToInt32($hex.Substring(0, 2), 16)
    $g = [Convert]::ToInt32($hex.Substring(2, 2), 16)
    $b = [Convert]::ToInt32($hex.Substring(4, 2), 16)
    $shape.Line.ForeColor.RGB = $r + ($g * 256) + ($b * 65536)
}

# Slide 1 ("p6" / title slide): drop the slide-number placeholder.
$s1 = $p.Slides.Item(1)
$s1.HeadersFooters.SlideNumber.Visible = $false

# Slide 4 ("p9"): recolor the two rounded-rectangle outlines.
$s4 = $p.Slides.Item(4)
Set-LineColorRGB $s4.Shapes.Item(4) "6F0A19"
Set-LineColorRGB $s4.Shapes.Item(5) "006778"

# Slide 8 ("p13" / closing slide): drop the slide-number placeholder.
$s8 = $p.Slides.Item(8)
$s8.HeadersFooters.SlideNumber.Visible = $false
